# Refresh the "auto-updating" date placeholder shown in the footer area of
# the slide master and every slide layout. PowerPoint recalculates this
# field (id {15F7C7F1-F004-4E77-BB99-7FE000EB8CD6}, type=datetimeFigureOut)
# against the current date whenever the deck is touched/saved; here the
# stored text moves from "12/7/2021" to "12/12/21".

$p = $ppt.ActivePresentation
$newDate = "12/12/21"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout under the master also carries its own copy of the
# date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
